$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.524.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.750.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4476"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07504"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.093"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.029"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.750.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001062"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06386"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.61%  "

$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("E21").Value = "  -2.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.859"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.564.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.083"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.949.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.091"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.082"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.658"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09023"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.555"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02298"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06018"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6361"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.944"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.206"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.384"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.770"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.46%  "

$ws.Range("E44").Value = "  -2.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.722"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5893"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.955"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.147"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06857"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.72%  "

